$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Version & History": add a new V1.10 change-history row
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Version & History")

# Duplicate the formatting of the last existing row (row 14) into the
# new row 15 so fonts/borders/number formats line up with the rest of
# the table.
$ws1.Range("A14:E14").Copy($ws1.Range("A15:E15"))

$ws1.Range("A15").Value() = "V1.10"

# Rich text cell: two runs with different fonts, same as the existing
# "V1.9" entry right above it.
$changeText = "Changed signals unit for the Radar Sensor:`n- Longitudinal RVX`n- Longitudinal EGO`n- Lateral RVY`n- Lateral EGO`n"
$ws1.Range("B15").Value() = $changeText
$headLen = 43
$headRun = $ws1.Range("B15").Characters(1, $headLen)
$headRun.Font.Name = "CorpoS"
$headRun.Font.Size = 11
$headRun.Font.Color = 0
$bodyRun = $ws1.Range("B15").Characters($headLen + 1, $changeText.Length - $headLen)
$bodyRun.Font.Name = "Calibri"
$bodyRun.Font.Size = 11
$bodyRun.Font.Color = 0

$ws1.Range("C15").Value() = "Zborai Attila"
$ws1.Range("D15").Value() = 42828
$ws1.Range("E15").Value() = "Draft version"

$ws1.Rows.Item(14).RowHeight = 76.1
$ws1.Rows.Item(15).RowHeight = 65.65

# ---------------------------------------------------------------
# Sheet "CommunicationMatrix": Radar Sensor signals now report
# their distance in on-screen Pixels instead of km/h or m.
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("CommunicationMatrix")
$ws2.Cells.Item(15, 6).Value() = "Pixel"  # Longitudinal RVX
$ws2.Cells.Item(16, 6).Value() = "Pixel"  # Longitudinal EGO
$ws2.Cells.Item(17, 6).Value() = "Pixel"  # Lateral RVY
$ws2.Cells.Item(18, 6).Value() = "Pixel"  # Lateral EGO

$ws2.Activate()
$ws2.Range("F22").Select()

# Re-activate the "Version & History" sheet last so it stays the tab
# that is selected when the workbook is reopened (matches the source).
$ws1.Activate()
$ws1.Range("B15").Select()
Write-Host "done"
